$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Drop the trailing rows (94:132) that fall out of the new 92-row window
$ws.Rows("94:132").Delete()

# Roll the weekly WTREGEN series forward: new dates/values for rows 2:93
$ws.Cells.Item(2,1).Value = 44475
$ws.Cells.Item(2,2).Value = 135.582
$ws.Cells.Item(3,1).Value = 44482
$ws.Cells.Item(3,2).Value = 77.858
$ws.Cells.Item(4,1).Value = 44489
$ws.Cells.Item(4,2).Value = 83.386
$ws.Cells.Item(5,1).Value = 44496
$ws.Cells.Item(5,2).Value = 213.291
$ws.Cells.Item(6,1).Value = 44503
$ws.Cells.Item(6,2).Value = 280.205
$ws.Cells.Item(7,1).Value = 44510
$ws.Cells.Item(7,2).Value = 259.293
$ws.Cells.Item(8,1).Value = 44517
$ws.Cells.Item(8,2).Value = 211.453
$ws.Cells.Item(9,1).Value = 44524
$ws.Cells.Item(9,2).Value = 164.087
$ws.Cells.Item(10,1).Value = 44531
$ws.Cells.Item(10,2).Value = 152.651
$ws.Cells.Item(11,1).Value = 44538
$ws.Cells.Item(11,2).Value = 115.117
$ws.Cells.Item(12,1).Value = 44545
$ws.Cells.Item(12,2).Value = 104.798
$ws.Cells.Item(13,1).Value = 44552
$ws.Cells.Item(13,2).Value = 84.954
$ws.Cells.Item(14,1).Value = 44559
$ws.Cells.Item(14,2).Value = 211.849
$ws.Cells.Item(15,1).Value = 44566
$ws.Cells.Item(15,2).Value = 389.751
$ws.Cells.Item(16,1).Value = 44573
$ws.Cells.Item(16,2).Value = 456.401
$ws.Cells.Item(17,1).Value = 44580
$ws.Cells.Item(17,2).Value = 488.697
$ws.Cells.Item(18,1).Value = 44587
$ws.Cells.Item(18,2).Value = 614.877
$ws.Cells.Item(19,1).Value = 44594
$ws.Cells.Item(19,2).Value = 675.773
$ws.Cells.Item(20,1).Value = 44601
$ws.Cells.Item(20,2).Value = 686.384
$ws.Cells.Item(21,1).Value = 44608
$ws.Cells.Item(21,2).Value = 692.541
$ws.Cells.Item(22,1).Value = 44615
$ws.Cells.Item(22,2).Value = 697.838
$ws.Cells.Item(23,1).Value = 44622
$ws.Cells.Item(23,2).Value = 682.903
$ws.Cells.Item(24,1).Value = 44629
$ws.Cells.Item(24,2).Value = 652.039
$ws.Cells.Item(25,1).Value = 44636
$ws.Cells.Item(25,2).Value = 570.606
$ws.Cells.Item(26,1).Value = 44643
$ws.Cells.Item(26,2).Value = 614.616
$ws.Cells.Item(27,1).Value = 44650
$ws.Cells.Item(27,2).Value = 575.065
$ws.Cells.Item(28,1).Value = 44657
$ws.Cells.Item(28,2).Value = 574.451
$ws.Cells.Item(29,1).Value = 44664
$ws.Cells.Item(29,2).Value = 547.308
$ws.Cells.Item(30,1).Value = 44671
$ws.Cells.Item(30,2).Value = 711.406
$ws.Cells.Item(31,1).Value = 44678
$ws.Cells.Item(31,2).Value = 944.328
$ws.Cells.Item(32,1).Value = 44685
$ws.Cells.Item(32,2).Value = 945.478
$ws.Cells.Item(33,1).Value = 44692
$ws.Cells.Item(33,2).Value = 945.658
$ws.Cells.Item(34,1).Value = 44699
$ws.Cells.Item(34,2).Value = 886.964
$ws.Cells.Item(35,1).Value = 44706
$ws.Cells.Item(35,2).Value = 821.535
$ws.Cells.Item(36,1).Value = 44713
$ws.Cells.Item(36,2).Value = 789.532
$ws.Cells.Item(37,1).Value = 44720
$ws.Cells.Item(37,2).Value = 723.384
$ws.Cells.Item(38,1).Value = 44727
$ws.Cells.Item(38,2).Value = 656.87
$ws.Cells.Item(39,1).Value = 44734
$ws.Cells.Item(39,2).Value = 758.283
$ws.Cells.Item(40,1).Value = 44741
$ws.Cells.Item(40,2).Value = 756.627
$ws.Cells.Item(41,1).Value = 44748
$ws.Cells.Item(41,2).Value = 700.457
$ws.Cells.Item(42,1).Value = 44755
$ws.Cells.Item(42,2).Value = 652.572
$ws.Cells.Item(43,1).Value = 44762
$ws.Cells.Item(43,2).Value = 613.878
$ws.Cells.Item(44,1).Value = 44769
$ws.Cells.Item(44,2).Value = 602.945
$ws.Cells.Item(45,1).Value = 44776
$ws.Cells.Item(45,2).Value = 594.115
$ws.Cells.Item(46,1).Value = 44783
$ws.Cells.Item(46,2).Value = 557.265
$ws.Cells.Item(47,1).Value = 44790
$ws.Cells.Item(47,2).Value = 545.321
$ws.Cells.Item(48,1).Value = 44797
$ws.Cells.Item(48,2).Value = 535.267
$ws.Cells.Item(49,1).Value = 44804
$ws.Cells.Item(49,2).Value = 612.536
$ws.Cells.Item(50,1).Value = 44811
$ws.Cells.Item(50,2).Value = 581.295
$ws.Cells.Item(51,1).Value = 44818
$ws.Cells.Item(51,2).Value = 593.808
$ws.Cells.Item(52,1).Value = 44825
$ws.Cells.Item(52,2).Value = 692.496
$ws.Cells.Item(53,1).Value = 44832
$ws.Cells.Item(53,2).Value = 689.569
$ws.Cells.Item(54,1).Value = 44839
$ws.Cells.Item(54,2).Value = 633.939
$ws.Cells.Item(55,1).Value = 44846
$ws.Cells.Item(55,2).Value = 608.302
$ws.Cells.Item(56,1).Value = 44853
$ws.Cells.Item(56,2).Value = 607.199
$ws.Cells.Item(57,1).Value = 44860
$ws.Cells.Item(57,2).Value = 634.548
$ws.Cells.Item(58,1).Value = 44867
$ws.Cells.Item(58,2).Value = 598.544
$ws.Cells.Item(59,1).Value = 44874
$ws.Cells.Item(59,2).Value = 527.479
$ws.Cells.Item(60,1).Value = 44881
$ws.Cells.Item(60,2).Value = 502.982
$ws.Cells.Item(61,1).Value = 44888
$ws.Cells.Item(61,2).Value = 479.474
$ws.Cells.Item(62,1).Value = 44895
$ws.Cells.Item(62,2).Value = 511.474
$ws.Cells.Item(63,1).Value = 44902
$ws.Cells.Item(63,2).Value = 432.335
$ws.Cells.Item(64,1).Value = 44909
$ws.Cells.Item(64,2).Value = 355.517
$ws.Cells.Item(65,1).Value = 44916
$ws.Cells.Item(65,2).Value = 459.78
$ws.Cells.Item(66,1).Value = 44923
$ws.Cells.Item(66,2).Value = 427.926
$ws.Cells.Item(67,1).Value = 44930
$ws.Cells.Item(67,2).Value = 423.625
$ws.Cells.Item(68,1).Value = 44937
$ws.Cells.Item(68,2).Value = 372.34
$ws.Cells.Item(69,1).Value = 44944
$ws.Cells.Item(69,2).Value = 339.018
$ws.Cells.Item(70,1).Value = 44951
$ws.Cells.Item(70,2).Value = 491.848
$ws.Cells.Item(71,1).Value = 44958
$ws.Cells.Item(71,2).Value = 560.089
$ws.Cells.Item(72,1).Value = 44965
$ws.Cells.Item(72,2).Value = 493.277
$ws.Cells.Item(73,1).Value = 44972
$ws.Cells.Item(73,2).Value = 490.379
$ws.Cells.Item(74,1).Value = 44979
$ws.Cells.Item(74,2).Value = 477.333
$ws.Cells.Item(75,1).Value = 44986
$ws.Cells.Item(75,2).Value = 381.245
$ws.Cells.Item(76,1).Value = 44993
$ws.Cells.Item(76,2).Value = 333.35
$ws.Cells.Item(77,1).Value = 45000
$ws.Cells.Item(77,2).Value = 232.866
$ws.Cells.Item(78,1).Value = 45007
$ws.Cells.Item(78,2).Value = 259.587
$ws.Cells.Item(79,1).Value = 45014
$ws.Cells.Item(79,2).Value = 183.577
$ws.Cells.Item(80,1).Value = 45021
$ws.Cells.Item(80,2).Value = 168.793
$ws.Cells.Item(81,1).Value = 45028
$ws.Cells.Item(81,2).Value = 109.208
$ws.Cells.Item(82,1).Value = 45035
$ws.Cells.Item(82,2).Value = 166.555
$ws.Cells.Item(83,1).Value = 45042
$ws.Cells.Item(83,2).Value = 291.702
$ws.Cells.Item(84,1).Value = 45049
$ws.Cells.Item(84,2).Value = 269.216
$ws.Cells.Item(85,1).Value = 45056
$ws.Cells.Item(85,2).Value = 197.666
$ws.Cells.Item(86,1).Value = 45063
$ws.Cells.Item(86,2).Value = 116.22
$ws.Cells.Item(87,1).Value = 45070
$ws.Cells.Item(87,2).Value = 61.952
$ws.Cells.Item(88,1).Value = 45077
$ws.Cells.Item(88,2).Value = 48.954
$ws.Cells.Item(89,1).Value = 45084
$ws.Cells.Item(89,2).Value = 44.756
$ws.Cells.Item(90,1).Value = 45091
$ws.Cells.Item(90,2).Value = 102.118
$ws.Cells.Item(91,1).Value = 45098
$ws.Cells.Item(91,2).Value = 276.85
$ws.Cells.Item(92,1).Value = 45105
$ws.Cells.Item(92,2).Value = 390.571
$ws.Cells.Item(93,1).Value = 45112
$ws.Cells.Item(93,2).Value = 415.441

$ws2 = $wb.Worksheets.Item("SeriesInfo")
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "2023-07-09"
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "2023-07-09"
$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("B7").Value = "2023-07-05"
$ws2.Range("B14").Value = "2023-07-06 15:34:05-05"
$ws2.Range("B15").Value = 84
